$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Regenerate merged AHB file headers -------------------------------
# The sheet compares two message-implementation-guide versions side by
# side: columns A-J describe the FV2410 variant, K holds the textual
# "diff" marker, and columns L-U describe the FV2504 variant. Rename the
# generic "_old" / "_new" header suffixes to the concrete version tags.
$headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410",
    "diff",
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2. Turn the data range into a real Excel Table ----------------------
# Adds xl/tables/table1.xml, the sheet1 <tableParts> reference and the
# worksheet-part relationship: column names are picked up from row 1.
$lastRow = $ws.UsedRange.Rows.Count
$lastCol = $ws.UsedRange.Columns.Count
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"

# --- 3. Freeze the header row --------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
